# Add a new leetcoder entry to the ranking sheet.
# This inserts a new row above the existing row 17 (shifting rows 17-34
# down to 18-35) and fills in the new row with the Edwards310 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 17 - everything below shifts down by one.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17.
$ws.Cells.Item(17, 1).Value = 32468
$ws.Cells.Item(17, 2).Value = "https://leetcode.com/u/Edwards310/"
$ws.Cells.Item(17, 3).Value = 437
$ws.Cells.Item(17, 4).Value = 483
$ws.Cells.Item(17, 5).Value = 16
$ws.Cells.Item(17, 12).Value = "No data"
$ws.Cells.Item(17, 14).Value = 406
$ws.Cells.Item(17, 15).Value = 336
$ws.Cells.Item(17, 16).Value = 68

# Match the number format used by the other "Rank" column cells (style s="1").
$ws.Cells.Item(17, 1).NumberFormat = "#,##0"

# Update the selected cell, matching the author's saved selection state.
$ws.Range("F17").Select()
